# "four scripts all plots - corrected graphics"
# Clear the stray/incorrect helper values in columns G:I (rows 2-7) on
# Sheet1 that were left over from an earlier graphics pass, while keeping
# the cell formatting (number styles) intact, then leave the new working
# selection on I2:I6 (matching where the data used to live).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the leftover numeric values in G2:H7, but preserve their styles.
$ws.Range("G2:H7").ClearContents() | Out-Null

# Remove the leftover "1" placeholder values in I2:I6 entirely (these cells
# had no explicit style, so clearing drops them from the sheet data).
$ws.Range("I2:I6").ClearContents() | Out-Null

# Update the active selection to reflect the corrected working range.
$ws.Range("I2:I6").Select() | Out-Null
